$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tickets updated")

# Clear the old contents and formatting entirely, then rebuild
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "_id"
$ws.Range("B1").Value = "internalField"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "__v"
$ws.Range("E1").Value = "readableField"
$ws.Range("F1").Value = "type"
$ws.Range("G1").Value = "checked"

# Row 2: email
$ws.Range("B2").Value = "email"
$ws.Range("C2").Value = "core"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "Email Id"
$ws.Range("F2").Value = "string"
$ws.Range("G2").Value = $true

# Row 3: nickname
$ws.Range("B3").Value = "nickname"
$ws.Range("C3").Value = "core"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "Nick Name"
$ws.Range("F3").Value = "string"
$ws.Range("G3").Value = $true

# Row 4: phoneNumber
$ws.Range("B4").Value = "phoneNumber"
$ws.Range("C4").Value = "core"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "Phone Number"
$ws.Range("F4").Value = "string"
$ws.Range("G4").Value = $true

# Row 5: phoneNumberPrefix
$ws.Range("B5").Value = "phoneNumberPrefix"
$ws.Range("C5").Value = "core"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = "Phone Number Prefix"
$ws.Range("F5").Value = "string"
$ws.Range("G5").Value = $true

# Row 6: amount
$ws.Range("B6").Value = "amount"
$ws.Range("C6").Value = "core"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = "Amount"
$ws.Range("F6").Value = "number"
$ws.Range("G6").Value = $true

# Row 7: followUp
$ws.Range("B7").Value = "followUp"
$ws.Range("C7").Value = "core"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = "FollowUp"
$ws.Range("F7").Value = "date"
$ws.Range("G7").Value = $true

# Row 8: agree
$ws.Range("B8").Value = "agree"
$ws.Range("C8").Value = "core"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = "Agree"
$ws.Range("F8").Value = "boolean"
$ws.Range("G8").Value = $false

# Row 9: status
$ws.Range("B9").Value = "status"
$ws.Range("C9").Value = "core"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = "Status"
$ws.Range("F9").Value = "string"
$ws.Range("G9").Value = $false
